$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.502.78'
$ws.Range('E2').Value = '  +4.28%  '
$ws.Range('D3').Value = '3.361.94'
$ws.Range('E3').Value = '  +9.23%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '256.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '623.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.47%  '
$ws.Range('E7').Value = '  +9.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.388'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.97%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = '3.360.32'
$ws.Range('E10').Value = '  +9.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.792'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('E12').Value = '  +1.77%  '
$ws.Range('D13').Value = '98.191.93'
$ws.Range('E13').Value = '  +4.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000246'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.97%  '
$ws.Range('D16').Value = '3.982.62'
$ws.Range('E16').Value = '  +9.24%  '
$ws.Range('E17').Value = '  +3.53%  '
$ws.Range('D18').Value = '3.360.46'
$ws.Range('E18').Value = '  +10.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '485.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.72%  '
$ws.Range('E23').Value = '  +10.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.71%  '
$ws.Range('E25').Value = '  +3.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = '3.537.97'
$ws.Range('E28').Value = '  +9.64%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.187'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.244'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.152'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '515.53'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.95'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.66%  '
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.25'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.784'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +17.34%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '160.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.99%  '
$ws.Range('E49').Value = '  +7.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.41'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.51'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.56%  '
